# baffin bay floats parameters
# Adds a new "column E" set of values to the ARVOR float parameter sheet,
# highlights the amended B-column cells in yellow, and tidies up the
# saved view state (top-left cell / selection) left over from editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Highlight the parameter rows that received a new value in column E.
$ws.Range("B3").Interior.Color = 65535
$ws.Range("B4").Interior.Color = 65535
$ws.Range("B5").Interior.Color = 65535
$ws.Range("B16").Interior.Color = 65535
$ws.Range("B62").Interior.Color = 65535

# New float-specific values in column E.
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 24
$ws.Range("E5").Value = 245
$ws.Range("E16").Value = 1000
$ws.Range("E62").Value = 0
$ws.Range("E3:E79").Font.Size = 10
